# Weekly update: insert 4 new observation rows (new rows 141-144) above the
# existing data, pushing the former rows 141-168 down to 145-172.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 141 (shifts old rows 141:168 -> 145:172).
$ws.Range("A141:A144").EntireRow.Insert()

# Common column values shared by every row in this data set.
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region  = "Metropolitana"
$codreg  = 13
$catId   = 100112026
$categoria = "Haba"
$variedad  = "Sin especificar"
$unidad    = "`$/saco 25 kilos"
$kgUnid    = 25
$clasificacion = "Hortaliza"

# New row 141
$r = 141
$ws.Cells.Item($r, 1).Value = 6
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44505
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 380
$ws.Cells.Item($r, 11).Value = 6000
$ws.Cells.Item($r, 12).Value = 7000
$ws.Cells.Item($r, 13).Value = 6395
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Región Metropolitana"
$ws.Cells.Item($r, 16).Value = 256
$ws.Cells.Item($r, 17).Value = $kgUnid
$ws.Cells.Item($r, 18).Value = $clasificacion

# New row 142
$r = 142
$ws.Cells.Item($r, 1).Value = 6
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44505
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 1410
$ws.Cells.Item($r, 11).Value = 6000
$ws.Cells.Item($r, 12).Value = 7000
$ws.Cells.Item($r, 13).Value = 6468
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Región del Maule"
$ws.Cells.Item($r, 16).Value = 259
$ws.Cells.Item($r, 17).Value = $kgUnid
$ws.Cells.Item($r, 18).Value = $clasificacion

# New row 143
$r = 143
$ws.Cells.Item($r, 1).Value = 6
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44505
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = "Segunda"
$ws.Cells.Item($r, 10).Value = 100
$ws.Cells.Item($r, 11).Value = 5000
$ws.Cells.Item($r, 12).Value = 5000
$ws.Cells.Item($r, 13).Value = 5000
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Región Metropolitana"
$ws.Cells.Item($r, 16).Value = 200
$ws.Cells.Item($r, 17).Value = $kgUnid
$ws.Cells.Item($r, 18).Value = $clasificacion

# New row 144
$r = 144
$ws.Cells.Item($r, 1).Value = 6
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44505
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = "Segunda"
$ws.Cells.Item($r, 10).Value = 550
$ws.Cells.Item($r, 11).Value = 5000
$ws.Cells.Item($r, 12).Value = 5000
$ws.Cells.Item($r, 13).Value = 5000
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Región del Maule"
$ws.Cells.Item($r, 16).Value = 200
$ws.Cells.Item($r, 17).Value = $kgUnid
$ws.Cells.Item($r, 18).Value = $clasificacion
